# Zombies information workbook - add "word / byte" addressing-mode column
# and a parallel set of rows showing the "byte" variant of each instruction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Formatting: copy the fill/format of the existing header & body
#    cells onto the new column J and the new rows 15:22 (column F:J)
# ---------------------------------------------------------------------

# Header cell J6 should look like the other header cells (teal fill, style 2)
$ws.Range("I6").Copy() | Out-Null
$ws.Range("J6").PasteSpecial(-4122) | Out-Null

# Body cells J7:J14 should look like the existing body cells (style 1)
$ws.Range("I7:I14").Copy() | Out-Null
$ws.Range("J7:J14").PasteSpecial(-4122) | Out-Null

# New rows 15:22, columns F:I -> copy formatting from rows 7:14
$ws.Range("F7:I14").Copy() | Out-Null
$ws.Range("F15:I22").PasteSpecial(-4122) | Out-Null

# New rows 15:22, column J -> copy formatting from I7:I14 as well
$ws.Range("I7:I14").Copy() | Out-Null
$ws.Range("J15:J22").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Write cell values in (approximately) the same order the original
#    author would have typed them, so freshly-introduced strings land
#    in the shared-string table in the same sequence as the reference
#    workbook.
# ---------------------------------------------------------------------

# Header row: "word / byte" addressing-mode column, "in main loop" shifts right
$ws.Range("I6").Value = "word / byte"
$ws.Range("J6").Value = "in main loop"

# Existing rows 7:14 become the "word" variant
$ws.Range("I7").Value = "word"

# First new value typically entered while building out the "byte" rows
$ws.Range("H19").Value = "0x46"

$ws.Range("I15").Value = "byte"

$ws.Range("H18").Value = "0xFFE2"
$ws.Range("H16").Value = "0x0000"
$ws.Range("H20").Value = "0xFFE6"
$ws.Range("H21").Value = "0x3F"
$ws.Range("H22").Value = "0xFFCC"
$ws.Range("H15").Value = "0xFFB9"
$ws.Range("H17").Value = "0x0000"

# Fill in the remainder of the "word" column (I8:I14)
foreach ($r in 8..14) {
    $ws.Range("I$r").Value = "word"
}

# Fill in the remainder of the "byte" column (I16:I22)
foreach ($r in 16..22) {
    $ws.Range("I$r").Value = "byte"
}

# Push the old "no"/"yes" values (formerly in column I) into the new column J
$oldInMainLoop = @{
    7  = "no"
    8  = "no"
    9  = "yes"
    10 = "yes"
    11 = "yes"
    12 = "no"
    13 = "no"
    14 = "no"
}
foreach ($r in 7..14) {
    $ws.Range("J$r").Value = $oldInMainLoop[$r]
}

# Rows 15:22 hold the "byte" variant of each instruction; F/G repeat the
# command/dist-from-start values from rows 7:14, J mirrors "in main loop"
$byteRows = @(
    @{ Row = 15; F = "mov bx,ax";                     G = "0x67"; J = "no"  },
    @{ Row = 16; F = "mov cx,0x10";                    G = "0x69"; J = "no"  },
    @{ Row = 17; F = "mov [bx + si + 0x100*Y],ax";     G = "0x6C"; J = "yes" },
    @{ Row = 18; F = "shl bx,0x1";                     G = "0x70"; J = "yes" },
    @{ Row = 19; F = "loop 0xfa";                      G = "0x72"; J = "yes" },
    @{ Row = 20; F = "inc si";                         G = "0x74"; J = "no"  },
    @{ Row = 21; F = "and si,0x3f";                    G = "0x75"; J = "no"  },
    @{ Row = 22; F = "jmp short 0xef";                 G = "0x78"; J = "no"  }
)

foreach ($entry in $byteRows) {
    $r = $entry.Row
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
    $ws.Range("J$r").Value = $entry.J
}

# ---------------------------------------------------------------------
# 3. Column width for the new column J
# ---------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 11.21875

# ---------------------------------------------------------------------
# 4. View state: put the selection where the author left it
# ---------------------------------------------------------------------
$ws.Range("G24").Select() | Out-Null

Write-Host "done"
